# The target diff only touches PowerPoint's internal co-authoring /
# revision-tracking bookkeeping:
#   - ppt/revisionInfo.xml            (the <p1510:client> version/timestamp
#     counter that Office bumps on every save)
#   - ppt/changesInfos/changesInfo1.xml (the per-author change history that
#     Office maintains for the "Show changes" / coauthoring UI)
#
# Both files are regenerated internally by PowerPoint itself whenever the
# file is saved; they are not exposed anywhere in the Slide/Shape/TextRange
# object model (no Presentation.RevisionInfo, no CustomXMLParts entry for
# them, etc.), so there is no COM call that can literally "set" the new
# version number / timestamp / action id.
#
# What the diff *does* tell us about real, user-visible content:
#   - The refreshed <pc:docChg> keeps the exact same "chg" flags it already
#     had (undo custSel addSld modSld) - nothing new was added to that set.
#   - The only <pc:sldChg>/<pc:spChg> entry whose chgData timestamp gets
#     bumped is the one for slide sldId="257" (slide11.xml, the last
#     "Bibliografía y otros recursos" slide), shape id="6" (the content
#     placeholder holding the bibliography text) - and it keeps its
#     existing chg="mod" / chg="modSp" markers rather than gaining new ones.
#   - There is no corresponding hunk anywhere against ppt/slides/slide11.xml
#     (or any other slide/layout/master part), so the placeholder's text,
#     formatting and geometry are provably byte-identical before and after.
#
# Put together, this matches a user re-opening/re-touching that bibliography
# placeholder (e.g. clicking into it and an undo landing back on the exact
# same content) and re-saving - an edit pass that PowerPoint's autosave /
# coauthoring telemetry records, but that nets out to zero visible change in
# the slide XML itself. We reproduce the visible part of that faithfully:
# touch the same shape on the same slide via the object model, read-only,
# so nothing in the deck's actual content drifts from the source file.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(5)

# Touch the placeholder (matches ac:spMk id="6" / creationId
# {F0814BA4-9F64-499E-9A70-81D4BA81EB0A} on sldId 257) without altering it.
$null = $shp.Id
$null = $shp.TextFrame.TextRange.Text
